# Scheduled runner update: refresh computed market-price / profit columns
# (H: currentAveragePrice, I: currentAveragePriceNQ, J: currentAveragePriceHQ,
#  K: LevePriceNQ, L: LevePriceHQ, M: LeveProfitNQ, N: LeveProfitHQ)
# for a batch of Leve rows across several job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 804.44446
$ws.Range("I33").Value = 569.6429000000001
$ws.Range("K33").Value = 569.6429000000001
$ws.Range("M33").Value = -340.6429000000001

$ws.Range("H98").Value = 4368.7646
$ws.Range("I98").Value = 3086.75
$ws.Range("J98").Value = 5508.3335
$ws.Range("K98").Value = 3086.75
$ws.Range("L98").Value = 5508.3335
$ws.Range("M98").Value = -1588.75
$ws.Range("N98").Value = -8504.333500000001

$ws.Range("H122").Value = 4368.7646
$ws.Range("I122").Value = 3086.75
$ws.Range("J122").Value = 5508.3335
$ws.Range("K122").Value = 9260.25
$ws.Range("L122").Value = 16525.0005
$ws.Range("M122").Value = -6810.25
$ws.Range("N122").Value = -21425.0005

$ws.Range("H135").Value = 768.7273
$ws.Range("I135").Value = 795.3684
$ws.Range("J135").Value = 600
$ws.Range("K135").Value = 7158.3156
$ws.Range("L135").Value = 5400
$ws.Range("M135").Value = -4623.3156
$ws.Range("N135").Value = -10470

$ws.Range("H137").Value = 1828.5483
$ws.Range("I137").Value = 1255.6
$ws.Range("J137").Value = 2365.6875
$ws.Range("K137").Value = 3766.8
$ws.Range("L137").Value = 7097.0625
$ws.Range("M137").Value = -1216.8
$ws.Range("N137").Value = -12197.0625

$ws.Range("H141").Value = 542891.9
$ws.Range("I141").Value = 1532.9412
$ws.Range("K141").Value = 4598.8236
$ws.Range("M141").Value = 581.1764000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5001568.5
$ws.Range("I2").Value = 11364447
$ws.Range("J2").Value = 2164.5
$ws.Range("K2").Value = 11364447
$ws.Range("L2").Value = 2164.5
$ws.Range("M2").Value = -11364334
$ws.Range("N2").Value = -2390.5

$ws.Range("H32").Value = 5106.8096
$ws.Range("I32").Value = 3621.8076
$ws.Range("K32").Value = 3621.8076
$ws.Range("M32").Value = -3334.8076

$ws.Range("H45").Value = 1579.6216
$ws.Range("I45").Value = 1035.4062
$ws.Range("J45").Value = 5062.6
$ws.Range("K45").Value = 1035.4062
$ws.Range("L45").Value = 5062.6
$ws.Range("M45").Value = -658.4061999999999
$ws.Range("N45").Value = -5816.6

$ws.Range("H61").Value = 2357.26
$ws.Range("I61").Value = 752.28
$ws.Range("J61").Value = 3962.24
$ws.Range("K61").Value = 752.28
$ws.Range("L61").Value = 3962.24
$ws.Range("M61").Value = -540.28
$ws.Range("N61").Value = -4386.24

$ws.Range("H74").Value = 772.9
$ws.Range("I74").Value = 434
$ws.Range("J74").Value = 1281.25
$ws.Range("K74").Value = 434
$ws.Range("L74").Value = 1281.25
$ws.Range("M74").Value = 440
$ws.Range("N74").Value = -3029.25

$ws.Range("H77").Value = 772.9
$ws.Range("I77").Value = 434
$ws.Range("J77").Value = 1281.25
$ws.Range("K77").Value = 2170
$ws.Range("L77").Value = 6406.25
$ws.Range("M77").Value = 2198
$ws.Range("N77").Value = -15142.25

$ws.Range("H102").Value = 6450
$ws.Range("I102").Value = 2900
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 2900
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = -1278
$ws.Range("N102").Value = -13244

$ws.Range("H116").Value = 5001568.5
$ws.Range("I116").Value = 11364447
$ws.Range("J116").Value = 2164.5
$ws.Range("K116").Value = 11364447
$ws.Range("L116").Value = 2164.5
$ws.Range("M116").Value = -11362153
$ws.Range("N116").Value = -6752.5

$ws.Range("H122").Value = 1608
$ws.Range("I122").Value = 1140.6666
$ws.Range("K122").Value = 3421.9998
$ws.Range("M122").Value = -971.9998000000001

$ws.Range("H132").Value = 18185128
$ws.Range("I132").Value = 27781126
$ws.Range("J132").Value = 3235.6843
$ws.Range("K132").Value = 83343378
$ws.Range("L132").Value = 9707.052899999999
$ws.Range("M132").Value = -83340848
$ws.Range("N132").Value = -14767.0529

$ws.Range("H136").Value = 2357.26
$ws.Range("I136").Value = 752.28
$ws.Range("J136").Value = 3962.24
$ws.Range("K136").Value = 2256.84
$ws.Range("L136").Value = 11886.72
$ws.Range("M136").Value = 293.1599999999999
$ws.Range("N136").Value = -16986.72

$ws.Range("H139").Value = 42500
$ws.Range("J139").Value = 42500
$ws.Range("L139").Value = 42500
$ws.Range("M139").Value = -52780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5001568.5
$ws.Range("I3").Value = 11364447
$ws.Range("J3").Value = 2164.5
$ws.Range("K3").Value = 11364447
$ws.Range("L3").Value = 2164.5
$ws.Range("M3").Value = -11364333
$ws.Range("N3").Value = -2392.5

$ws.Range("H105").Value = 1880.9584
$ws.Range("I105").Value = 1465.7142
$ws.Range("J105").Value = 2462.3
$ws.Range("K105").Value = 1465.7142
$ws.Range("L105").Value = 2462.3
$ws.Range("M105").Value = 281.2858000000001
$ws.Range("N105").Value = -5956.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2113.8784
$ws.Range("I31").Value = 1302.2391
$ws.Range("J31").Value = 3447.2856
$ws.Range("K31").Value = 1302.2391
$ws.Range("L31").Value = 3447.2856
$ws.Range("M31").Value = -1007.2391
$ws.Range("N31").Value = -4037.2856

$ws.Range("H34").Value = 2113.8784
$ws.Range("I34").Value = 1302.2391
$ws.Range("J34").Value = 3447.2856
$ws.Range("K34").Value = 1302.2391
$ws.Range("L34").Value = 3447.2856
$ws.Range("M34").Value = -1100.2391
$ws.Range("N34").Value = -3851.2856

$ws.Range("H122").Value = 2823.923
$ws.Range("I122").Value = 1968
$ws.Range("J122").Value = 4749.75
$ws.Range("K122").Value = 5904
$ws.Range("L122").Value = 14249.25
$ws.Range("M122").Value = -3454
$ws.Range("N122").Value = -19149.25

$ws.Range("H132").Value = 3354
$ws.Range("I132").Value = 1865.2
$ws.Range("J132").Value = 5384.1816
$ws.Range("K132").Value = 5595.6
$ws.Range("L132").Value = 16152.5448
$ws.Range("M132").Value = -3065.6
$ws.Range("N132").Value = -21212.5448

$ws.Range("H134").Value = 1603.9697
$ws.Range("I134").Value = 1256.4423
$ws.Range("J134").Value = 2894.7856
$ws.Range("K134").Value = 3769.3269
$ws.Range("L134").Value = 8684.356800000001
$ws.Range("M134").Value = -1234.3269
$ws.Range("N134").Value = -13754.3568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 246.15384
$ws.Range("J12").Value = 262.5
$ws.Range("L12").Value = 787.5
$ws.Range("N12").Value = -1133.5

$ws.Range("H131").Value = 1341.6562
$ws.Range("I131").Value = 1293.0769
$ws.Range("J131").Value = 1374.8948
$ws.Range("K131").Value = 3879.2307
$ws.Range("L131").Value = 4124.6844
$ws.Range("M131").Value = 1160.7693
$ws.Range("N131").Value = -14204.6844

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4541.231
$ws.Range("I132").Value = 3160.5715
$ws.Range("J132").Value = 6152
$ws.Range("K132").Value = 9481.7145
$ws.Range("L132").Value = 18456
$ws.Range("M132").Value = -6951.7145
$ws.Range("N132").Value = -23516

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2271.4285
$ws.Range("I68").Value = 1072.7273
$ws.Range("J68").Value = 6666.6665
$ws.Range("K68").Value = 1072.7273
$ws.Range("L68").Value = 6666.6665
$ws.Range("M68").Value = -323.7273
$ws.Range("N68").Value = -8164.6665

$ws.Range("H71").Value = 2271.4285
$ws.Range("I71").Value = 1072.7273
$ws.Range("J71").Value = 6666.6665
$ws.Range("K71").Value = 5363.636500000001
$ws.Range("L71").Value = 33333.3325
$ws.Range("M71").Value = -1619.636500000001
$ws.Range("N71").Value = -40821.3325
